# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" and before "2021-Q4",
#    with the same per-fund detail layout as the other quarter sheets.
# 2. Insert a new summary row at the top of the "总计" sheet's data table for
#    the 2022-Q3 totals, pushing the existing quarters down by one row.
# 3. Restore the original active sheet ("2021-Q1", the last tab).

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q4old = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q3" sheet, positioned after "总计" -------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Re-fetch a fresh handle to the 2021-Q4 sheet to use as a formatting
# template (stale references can silently repoint after a sheet insert).
$q4old = $wb.Worksheets.Item("2021-Q4")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Copy header row formatting (borders/bold/center/style) from the 2021-Q4
# sheet so the new sheet matches the existing look. Column A only carries a
# style on row 2 (the "index" column), so copy that separately.
$q4old.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q4old.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Data row - B,D,E,F,G are text-like values (leading zeros / fixed decimals
# must survive), so force text storage before assigning them.
$q3.Range("A2").Value = 0
$q3.Range("B2:G2").NumberFormat = "@"
$q3.Range("B2").Value = "001704"
$q3.Range("C2").Value = "国投瑞银进宝灵活配置混合"
$q3.Range("D2").Value = "42.10"
$q3.Range("E2").Value = "91.58"
$q3.Range("F2").Value = "5.03"
$q3.Range("G2").Value = "2.1176"
$q3.Range("H2").Value = 9

# --- 2. Insert the 2022-Q3 summary row into "总计" --------------------------
$total = $wb.Worksheets.Item("总计")

# Shift the existing summary rows down one (bottom-up so we don't clobber).
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.24
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.49

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.84

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.19

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 2.12

# --- 3. Restore the original active tab -------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
